$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 46044, 12846.0615987304, 12621.3137820134, 20883.86, 8136.07113246202, -5.26979523019236),
    @(3, 46045, 12800.3679595451, 11729.5385078364, 12075.86, 8663.57432832953, 346.552201506915),
    @(4, 46046, 4904.18741062795, 7777.38675817716, 12075.86, 8016.21838489182, 154.906047627874),
    @(5, 46047, 4719.31680199932, 7694.78360714886, 12075.86, 7902.11596425204, 146.709982141704),
    @(6, 46048, 12010.2628603363, 11681.5402910665, 12075.86, 8057.5523674214, 319.301360770329),
    @(7, 46049, 12315.9682835607, 12138.825866146, 12075.86, 8314.33628491134, 349.054256294057),
    @(8, 46050, 12315.9682835607, 11439.7837695994, 12075.86, 8314.33628491134, 319.927502271282),
    @(9, 46051, 12315.9682835607, 11646.782596875, 12075.86, 8314.33628491134, 328.552453407763),
    @(10, 46052, 12315.9682835607, 11065.9491895018, 12075.86, 8314.33628491134, 304.351061433879),
    @(11, 46053, 4867.38022112383, 7784.3049024531, 12075.86, 7930.27449985425, 151.613308429473),
    @(12, 46054, 5046.61366744637, 7536.64654050504, 9503.86, 7890.32895980472, 246.796479179573),
    @(13, 46055, 11445.6138712783, 10864.9981164976, 9503.86, 7806.04182518022, 381.965830903243),
    @(14, 46056, 11445.6138712783, 10814.8345414041, 9503.86, 7806.04182518022, 379.875681941013),
    @(15, 46057, 11445.6138712783, 10692.1617347382, 9503.86, 7806.04182518022, 374.7643149966)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
